# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns for the
# a816aaa2-1c96-441b-a7b2-68de5b93d7d3 row on both the zh-cn and de-de
# sheets, now that a (stale) handback has been received for that file.

$wb = $excel.ActiveWorkbook

$fileMd  = "a816aaa2-1c96-441b-a7b2-68de5b93d7d3.md"
$fileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b74f1ae55de759b21a2a281e863c5e5670a8a9a/e2e/a816aaa2-1c96-441b-a7b2-68de5b93d7d3.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e545fe6b34c2828d2b3c90a88c0fc7abe7fc5b7/e2e/a816aaa2-1c96-441b-a7b2-68de5b93d7d3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b74f1ae55de759b21a2a281e863c5e5670a8a9a/e2e/a816aaa2-1c96-441b-a7b2-68de5b93d7d3.md."

# --- zh-cn sheet -----------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I8").Value = $fileMd
$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $fileUrl, [Type]::Missing, [Type]::Missing, $fileMd)
$wsZh.Range("J8").Value = "a816aaa2-1c96-441b-a7b2-68de5b93d7d3.ca2cf66444c8301d4708f797a59b05fdfa94235b.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-09-01 04:47:30"
$wsZh.Range("P8").Value = $errorDetail

# widen the Error Detail column so the message is readable
$wsZh.Range("P1").ColumnWidth = $wsZh.Range("A1").ColumnWidth

# --- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I8").Value = $fileMd
$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $fileUrl, [Type]::Missing, [Type]::Missing, $fileMd)
$wsDe.Range("J8").Value = "a816aaa2-1c96-441b-a7b2-68de5b93d7d3.ca2cf66444c8301d4708f797a59b05fdfa94235b.de-de.xlf"
$wsDe.Range("K8").Value = "2016-09-01 04:47:38"
$wsDe.Range("P8").Value = $errorDetail

# widen the Error Detail column so the message is readable
$wsDe.Range("P1").ColumnWidth = $wsDe.Range("A1").ColumnWidth
